# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45557 to 45558 (one day later).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45557) {
        $cell.Value2 = 45558
    }
}
